$d = $word.ActiveDocument

# The texts for the ten new paragraphs to be inserted right after the
# "From Spec:" paragraph (these used to be squashed into one big paragraph).
$texts = @(
  '“There''s no such thing as perfection. You''re never finished with a film. You run out of time.” -- Peter Jackson, director of `The Lord of the Rings’ and ‘The Hobbit’ trilogies ',
  'One of the more difficult parts of project planning and execution is to define the scope and limits of the project. ',
  'As mentioned above, you never really complete project like these; all you can ever do is your best in the time available. ',
  'Part of that involves setting priorities and accepting that there will be features that will take too long to develop. ',
  'This means that it is important to set a scope for your project, as a means of ensuring that you make the most of the time available. ',
  'For example, if you are developing a game, you might consider only producing one level and two or three characters, in order to show a proof-ofconcept, rather than develop three levels and ten characters. ',
  'The scope is probably the most crucial part of your plan, and also the most difficult to define. One way to define the scope is to think of the deliverables for your project, i.e. what outcomes would you be able to show to someone who asks you to see the results of your work. ',
  'This will also include several statements about what will not be part of the project. For example, if you are using Open Street Maps to show the location of all your favourite shops, the deliverables would include the updated map, but not the Open Street Maps technology itself. ',
  'It would also not include many other features of Open Street Maps, or other interesting location -- just those which show your favourite shops. Also, be aware of the phenomenon of `scope creep'', which is the tendency for projects to incorporate more and more features. ',
  'There is nothing wrong with being ambitious, but you only have a certain amount of time. '
)

# Find the "From Spec:" paragraph -- new paragraphs get inserted right after it.
$anchor = $d.Paragraphs(6)

$idx = $anchor.Index
foreach ($t in $texts) {
    $cur = $d.Paragraphs($idx)
    $cur.Range.InsertParagraphAfter()
    $idx = $idx + 1
    $newPara = $d.Paragraphs($idx)
    $newPara.Range.Text = $t
}

# The big paragraph (formerly holding the whole text) now sits right after
# the newly-inserted paragraphs; trim it down to its final sentence only,
# keeping its paragraph mark formatting (italic) and the _GoBack bookmark.
$idx = $idx + 1
$big = $d.Paragraphs($idx)
$fullRange = $d.Range($big.Range.Start, $big.Range.End)
$fullRange.Text = 'At least one paragraph is expected.'

# Re-anchor the _GoBack bookmark to the start of that paragraph so it
# precedes the run, matching the original document order.
$trimmed = $d.Paragraphs($idx)
$bmRange = $d.Range($trimmed.Range.Start, $trimmed.Range.Start)
$d.Bookmarks.Add("_GoBack", $bmRange)
